# Amendments to xlsx and test summaries.
#
# The sheet originally held a single row: A1="a", B1="b", C1="c".
# It is reshaped into a 3-row x 2-column block (A1:B3) containing
# a,b / c,d / e,f (row-major), so C1's old value "c" is cleared and
# two more rows are appended with the new "c","d","e","f" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old C1 value - the data now lives in a 2-column range instead.
$ws.Range("C1").ClearContents()

# Row 2: c, d
$ws.Range("A2").Value = "c"
$ws.Range("B2").Value = "d"

# Row 3: e, f
$ws.Range("A3").Value = "e"
$ws.Range("B3").Value = "f"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("C6").Select()
